$wb = $excel.ActiveWorkbook

# Add the new sheet at the end of the workbook and name it "N=200000"
$newSheet = $wb.Worksheets.Add()
$newSheet.Name = "N=200000"
$lastIndex = $wb.Worksheets.Count
$newSheet.Move($null, $wb.Worksheets.Item($lastIndex))

$ws = $wb.Worksheets.Item("N=200000")

# Header row
$ws.Range("A1").Value = "Execução"
$ws.Range("B1").Value = "Tempo (ms)"

# Data rows
$ws.Range("A2").Value = 1
$ws.Range("B2").Value = "82.3002 ms"

$ws.Range("A3").Value = 2
$ws.Range("B3").Value = "83.2818 ms"

$ws.Range("A4").Value = 3
$ws.Range("B4").Value = "84.3840 ms"

$ws.Range("A5").Value = 4
$ws.Range("B5").Value = "84.4238 ms"

$ws.Range("A6").Value = 5
$ws.Range("B6").Value = "94.3730 ms"

$ws.Range("A7").Value = "Média"
$ws.Range("B7").Value = "85.7525 ms"

$ws.Range("A8").Value = "Desvio Padrão"
$ws.Range("B8").Value = "4.8983 ms"

$ws.Range("A1").Select()
